$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the last existing header cell (G1) onto
# the new header cell (H1) so it reuses the same header style, then set
# its value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell for the "Save" column.
$ws.Range("H2").Value = 0
